$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D8 text from "semente" (singular) to "sementes" (plural)
$ws.Range("D8").Value = "Raiz, caule, folhas, estróbilos e sementes"

# Update the active selection to F6 (as captured in the saved view state)
$ws.Range("F6").Select()
